# Montana overview workbook update.
#
# The source data export switched every numeric-looking metric cell over to
# plain text (inline string) values. On top of that, sheet "County" gained a
# new statewide "Total" row (55) and the eight previously-zeroed counties
# (rows 47-54) had their blank "0" placeholders replaced with formatted
# percent / currency text ("0.00%", "$0").
#
# Helper: force a cell to hold a literal text value (never auto-coerced back
# to a number by Excel), then drop the cell back to the default "Normal"
# style so we don't leave a stray number-format style behind.
function Set-TextValue($cell, $val) {
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overall": A2 (667) numeric -> text
# ---------------------------------------------------------------------
$wsOverall = $wb.Worksheets.Item("Overall")
Set-TextValue $wsOverall.Cells.Item(2,1) "667"

# ---------------------------------------------------------------------
# Sheet "County": B2:B46 numeric -> text (values unchanged)
# ---------------------------------------------------------------------
$wsCounty = $wb.Worksheets.Item("County")
$countyCounts = @(6, 7, 4, 12, 2, 41, 3, 5, 1, 3, 6, 2, 10, 54, 90, 9, 11, 4, 14, 58, 1, 15, 6, 1, 2, 1, 102, 2, 21, 2, 5, 3, 5, 17, 9, 3, 4, 9, 1, 24, 3, 1, 6, 1, 81)
$row = 2
foreach ($v in $countyCounts) {
    Set-TextValue $wsCounty.Cells.Item($row, 2) "$v"
    $row = $row + 1
}

# Rows 47-54: the eight zero-activity counties get reformatted placeholder
# text instead of bare "0" values.
for ($r = 47; $r -le 54; $r++) {
    Set-TextValue $wsCounty.Cells.Item($r, 2) "0.00%"
    Set-TextValue $wsCounty.Cells.Item($r, 3) "`$0"
    Set-TextValue $wsCounty.Cells.Item($r, 4) "0.00%"
    Set-TextValue $wsCounty.Cells.Item($r, 5) "0.00%"
    Set-TextValue $wsCounty.Cells.Item($r, 6) "0.00%"
}

# New row 55: statewide total, matching the other sheets' "Total" rows.
Set-TextValue $wsCounty.Cells.Item(55, 1) "Total"
Set-TextValue $wsCounty.Cells.Item(55, 2) "667"
Set-TextValue $wsCounty.Cells.Item(55, 3) "`$680,827,015"
Set-TextValue $wsCounty.Cells.Item(55, 4) "9.74%"
Set-TextValue $wsCounty.Cells.Item(55, 5) "-10.06%"
Set-TextValue $wsCounty.Cells.Item(55, 6) "63.72%"

# ---------------------------------------------------------------------
# Sheet "Congressional District": B2:B4 numeric -> text
# ---------------------------------------------------------------------
$wsDistrict = $wb.Worksheets.Item("Congressional District")
$districtCounts = @(358, 309, 667)
$row = 2
foreach ($v in $districtCounts) {
    Set-TextValue $wsDistrict.Cells.Item($row, 2) "$v"
    $row = $row + 1
}

# ---------------------------------------------------------------------
# Sheet "Size": B2:B8 numeric -> text
# ---------------------------------------------------------------------
$wsSize = $wb.Worksheets.Item("Size")
$sizeCounts = @(236, 189, 112, 42, 67, 21, 667)
$row = 2
foreach ($v in $sizeCounts) {
    Set-TextValue $wsSize.Cells.Item($row, 2) "$v"
    $row = $row + 1
}

# ---------------------------------------------------------------------
# Sheet "Subsector": B2:B13 numeric -> text
# ---------------------------------------------------------------------
$wsSubsector = $wb.Worksheets.Item("Subsector")
$subsectorCounts = @(45, 43, 66, 63, 28, 208, 1, 65, 8, 135, 5, 667)
$row = 2
foreach ($v in $subsectorCounts) {
    Set-TextValue $wsSubsector.Cells.Item($row, 2) "$v"
    $row = $row + 1
}
